# Actualización automática del index.html y archivo Excel
# The row for "Caso -496" (Ricardo Balbin 3851, OT 807846856) was removed.
# Deleting the entire row 50 shifts the subsequent rows (51 -> 50, 52 -> 51)
# up by one and shrinks the used range from N52 to N51, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Delete()
